$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top, shifting all existing content down by 2 rows
$ws.Rows("1:2").Insert()

# Populate the three new header rows with the citation / provenance notes
$ws.Range("A1").Value2 = "Data from Cameron EK, Cahill JF Jr, Bayne EM (2014) Root Foraging Influences Plant Growth Responses to Earthworm Foraging. PLoS ONE 9(9): e108873. doi:10.1371/journal.pone.0108873"
$ws.Range("A2").Value2 = "Original data is available at https://era.library.ualberta.ca/files/z029p5988#.WIKrJZJVeAA"
$ws.Range("A3").Value2 = "The data in this spreadsheet has had flaws introduced for educational purposes."

# Update the view so the new top row is selected and visible
$ws.Range("A2").Select()
